$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Estudiantes: " -> "Estudiante: " (collapses the split "Estudiante"/"s"/": "
#    runs into a single run, matching the target OOXML).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Estudiantes: ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Estudiante: ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "05 de Mayo del 2012" -> "Mayo del 2012"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("05 de Mayo del 2012", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Mayo del 2012", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "1ra. Corrección del diseño de la base de datos" paragraph ->
#    "2da. Corrección del diseño de la base de datos." split across four runs.
#    A straight Find/Replace would coalesce the whole sentence into a single
#    run, so the whole paragraph is rebuilt via InsertXML to keep the exact
#    run layout from the target document.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("1ra. Corrección del diseño de la base de datos") | Out-Null
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="008B7C80" w:rsidRDefault="001145E5" w:rsidP="00F441E4"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="0"/><w:rPr><w:lang w:val="es-CR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>2da</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="008B7C80"><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>Corrección del diseño de la base de datos</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng3.InsertXML($xml3) | Out-Null

# ---------------------------------------------------------------------------
# 4) Remove the trailing "Se dejará pendiente para una próxima iteración el
#    avatar." sentence (two runs) after "...]." for the Perfil de Usuario row.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" Se dejará pendiente para una próxima iteración el avatar.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Extend the "Lista de los 10 últimos libros..." sentence with the new
#    moderation / note text, split across many runs (plus proofErr markers),
#    matching the target document exactly - rebuilt via InsertXML.
# ---------------------------------------------------------------------------
$rng5 = $d.Content
$rng5.Find.Execute(" Lista de los 10 últimos libros publicados mostrados en la página principal.") | Out-Null
$xml5 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="002853F8" w:rsidRDefault="00577FE9" w:rsidP="004B74E5"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="0"/><w:rPr><w:lang w:val="es-CR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve"> [Código: </w:t></w:r><w:r w:rsidRPr="001145E5"><w:rPr><w:b/><w:lang w:val="es-CR"/></w:rPr><w:t>100</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="es-CR"/></w:rPr><w:t>9 Caso 1 Libro</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>].</w:t></w:r><w:r w:rsidR="004B74E5"><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve"> Lista de los 10 últimos libros publicados mostrados en la página principal</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve">, previa aprobación por el moderador de </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>contenidos</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>admin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>pwd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>: 123456)</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve">Nota.- </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>Falta la opción rechazar</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve"> contenido</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>, así como la categorización automática</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CR"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng5.InsertXML($xml5) | Out-Null

Write-Output "Edits applied."
